$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:E -> B:F) right
$ws.Range("A1").EntireColumn.Insert()

# New header cell value
$ws.Range("A1").Value = "FILE"

# Copy the header formatting (bold/centered style) from the neighboring header cell
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# New data cells for rows 2 and 3 (source file path)
$ws.Range("A2").Value = "sample\sampleSQL.xml"
$ws.Range("A3").Value = "sample\sampleSQL.xml"
